$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (A and B) to make room for data/loja, shifting existing columns right
$ws.Columns("A:B").Insert()

# Copy header style (from former A1, now C1) onto the two new header cells
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$ws.Range("A1").Value = 'data'
$ws.Range("B1").Value = 'loja'

# Fill new data/loja columns for all data rows
$ws.Range("A2:A10").Value = '30/07/2024'
$ws.Range("B2:B10").Value = 'sita9289590'

# Rewrite remaining columns (C..I) row by row to match final content,
# including the row5/row6 and row7/row8 swaps and updated link tracking_ids

# Row 2
$ws.Range("C2").Value = 'Controle Longa Distância Jfa Acqua 1200 Resistente A Água'
$ws.Range("D2").Value = 'ACQUA'
$ws.Range("E2").Value = 63.9
$ws.Range("F2").Value = 'Baixo'
$ws.Range("G2").Value = 'NA'
$ws.Range("H2").Value = 'classico'
$ws.Range("I2").Value = 'https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:154235731#searchVariation=MLB27687422&position=2&search_layout=stack&type=product&tracking_id=3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 3
$ws.Range("C3").Value = 'Conversor Fio Para Rca Remoto Slim 12v Jfa Automotivo Cd Dvd'
$ws.Range("D3").Value = 'Sem Modelo'
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = ''
$ws.Range("G3").Value = 'NA'
$ws.Range("H3").Value = 'classico'
$ws.Range("I3").Value = 'https://www.mercadolivre.com.br/conversor-fio-para-rca-remoto-slim-12v-jfa-automotivo-cd-dvd/p/MLB25707531?pdp_filters=seller_id:154235731#searchVariation=MLB25707531&position=4&search_layout=stack&type=product&tracking_id=3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 4
$ws.Range("C4").Value = 'Controle Remoto Universal Longa Distância Jfa K1200 Preto'
$ws.Range("D4").Value = 'K1200'
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 'Acima'
$ws.Range("G4").Value = 'NA'
$ws.Range("H4").Value = 'classico'
$ws.Range("I4").Value = 'https://www.mercadolivre.com.br/controle-remoto-universal-longa-distncia-jfa-k1200-preto/p/MLB28687615?pdp_filters=seller_id:154235731#searchVariation=MLB28687615&position=1&search_layout=stack&type=product&tracking_id=3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 5
$ws.Range("C5").Value = 'Filtro Anti-ruido Jfa Com Blindagem Eletromagnética 20/20k'
$ws.Range("D5").Value = 'Sem Modelo'
$ws.Range("E5").Value = 52.9
$ws.Range("F5").Value = ''
$ws.Range("G5").Value = 'NA'
$ws.Range("H5").Value = 'classico'
$ws.Range("I5").Value = 'https://produto.mercadolivre.com.br/MLB-4531110844-filtro-anti-ruido-jfa-com-blindagem-eletromagnetica-2020k-_JM#position%3D5%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 6
$ws.Range("C6").Value = 'Voltímetro Jfa Vs5hi 3 Em 1 Sequenciador High Voltagem /12v'
$ws.Range("D6").Value = 'Sem Modelo'
$ws.Range("E6").Value = 52.99
$ws.Range("F6").Value = ''
$ws.Range("G6").Value = 'NA'
$ws.Range("H6").Value = 'classico'
$ws.Range("I6").Value = 'https://produto.mercadolivre.com.br/MLB-4531096344-voltimetro-jfa-vs5hi-3-em-1-sequenciador-high-voltagem-12v-_JM#position%3D6%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 7
$ws.Range("C7").Value = 'Fonte Automotiva Jfa Bob Storm 90a Bivolt Carregador'
$ws.Range("D7").Value = 'FONTE 90 BOB'
$ws.Range("E7").Value = 435
$ws.Range("F7").Value = 'Acima'
$ws.Range("G7").Value = 'NA'
$ws.Range("H7").Value = 'classico'
$ws.Range("I7").Value = 'https://produto.mercadolivre.com.br/MLB-3629903553-fonte-automotiva-jfa-bob-storm-90a-bivolt-carregador-_JM#position%3D7%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 8
$ws.Range("C8").Value = 'Fonte Digital Jfa 70a Storm Carregador Inteligente Bateria '
$ws.Range("D8").Value = 'FONTE 70A STORM'
$ws.Range("E8").Value = 493.99
$ws.Range("F8").Value = 'Acima'
$ws.Range("G8").Value = 'NA'
$ws.Range("H8").Value = 'classico'
$ws.Range("I8").Value = 'https://produto.mercadolivre.com.br/MLB-3629883283-fonte-digital-jfa-70a-storm-carregador-inteligente-bateria-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 9
$ws.Range("C9").Value = 'Fonte Carregador De Bateria Jfa 40a Storm Som Automotivo'
$ws.Range("D9").Value = 'FONTE 40A STORM'
$ws.Range("E9").Value = 404
$ws.Range("F9").Value = 'Acima'
$ws.Range("G9").Value = 'NA'
$ws.Range("H9").Value = 'classico'
$ws.Range("I9").Value = 'https://produto.mercadolivre.com.br/MLB-3629872501-fonte-carregador-de-bateria-jfa-40a-storm-som-automotivo-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'

# Row 10
$ws.Range("C10").Value = 'Fonte Carregador Automotivo Jfa Storm 60a Bivolt Volt/amp'
$ws.Range("D10").Value = 'FONTE 60A STORM'
$ws.Range("E10").Value = 443.99
$ws.Range("F10").Value = 'Acima'
$ws.Range("G10").Value = 'NA'
$ws.Range("H10").Value = 'classico'
$ws.Range("I10").Value = 'https://produto.mercadolivre.com.br/MLB-3629847295-fonte-carregador-automotivo-jfa-storm-60a-bivolt-voltamp-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D3de79c9a-68f3-4eaf-9f21-b164c6734d5d'
